# Update the "+ MILIONÁRIA" lottery results sheet with the four newest
# draws (concursos 308-311), matching the data refresh performed in the
# source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Concurso, Bola1..Bola6, Trevo1, Trevo2 for draws 308-311 (rows 309-312)
$data = @(
    @(308, 12, 16, 18, 19, 35, 40, 1, 6),
    @(309, 1,  8,  15, 30, 37, 47, 2, 6),
    @(310, 15, 16, 33, 34, 40, 42, 2, 5),
    @(311, 19, 29, 32, 38, 42, 50, 2, 4)
)

$r = 309
foreach ($rowValues in $data) {
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
    $r++
}

# The previous four rows (305:308, draws 304-307) had been highlighted as
# the "latest" results - now that newer draws exist, drop that formatting
# back to the sheet's normal (unstyled) look ...
$ws.Range("A305:I308").ClearFormats()

# ... and mark the newly-added rows (309:312) as the latest results instead.
$ws.Range("A309:I312").NumberFormat = "General"

# Move the active selection off the old range and past the new rows, as
# happens after entering the last new value.
$ws.Range("J315").Select()
